# Unfinished Pinj calc function
# Adds a "Theta (rad)" column to the Steady State Bus Data sheet, computed
# from the existing Theta (now labeled "Theta (deg)") column, and leaves
# that sheet as the active tab/selection.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Steady State Bus Data")

# Insert a new column D (shifts old P/Q columns from D/E to E/F) to hold
# the bus angle in radians, next to the existing angle-in-degrees column.
$ws2.Columns(4).Insert()

# Relabel the existing angle column and label the new one.
$ws2.Range("C1").Value = "Theta (deg)"
$ws2.Range("D1").Value = "Theta (rad)"

# Pinj needs Theta in radians -- convert degrees to radians for each bus row.
$ws2.Range("D5").Formula = "=C5*PI()/180"
$ws2.Range("D6:D9").Formula = "=C6*PI()/180"

# Leave the bus data sheet as the active/selected sheet and cell.
$ws2.Activate() | Out-Null
$ws2.Range("L10").Select() | Out-Null
